# Case_3_123 (380 kV) line active-power results: pl_mw.xlsx, Sheet1
# Overwrite the simulated line-flow values in B2:O25 (cols A,E,H,L are
# the fixed index/zero columns and are left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 index 0)
$ws.Range("B2").Value = 0.3021441429092135
$ws.Range("C2").Value = 0.04032442125382829
$ws.Range("D2").Value = 0.2619894535071836
$ws.Range("F2").Value = 1.461030158445865
$ws.Range("G2").Value = 0.002461779634844103
$ws.Range("I2").Value = 0.6957237433609365
$ws.Range("J2").Value = 0.2981952928391962
$ws.Range("K2").Value = 0.3164671445548777
$ws.Range("M2").Value = 0.2868413457126735
$ws.Range("N2").Value = 1.739539872600418
$ws.Range("O2").Value = 3.33953822088418
# Row 3 (A3 index 1)
$ws.Range("B3").Value = 0.2715625377628612
$ws.Range("C3").Value = 0.03551205972293303
$ws.Range("D3").Value = 0.2569519675557501
$ws.Range("F3").Value = 1.463775374681902
$ws.Range("G3").Value = 0.002464017470458139
$ws.Range("I3").Value = 0.7009732478823452
$ws.Range("J3").Value = 0.2964387989239086
$ws.Range("K3").Value = 0.2830370935275255
$ws.Range("M3").Value = 0.2738324866643538
$ws.Range("N3").Value = 1.756508387189941
$ws.Range("O3").Value = 3.355606042519867
# Row 4 (A4 index 2)
$ws.Range("B4").Value = 0.2528213793436009
$ws.Range("C4").Value = 0.03255446131515782
$ws.Range("D4").Value = 0.2539700804974814
$ws.Range("F4").Value = 1.466158955167352
$ws.Range("G4").Value = 0.002465465875382556
$ws.Range("I4").Value = 0.7045127818257839
$ws.Range("J4").Value = 0.2955083224214832
$ws.Range("K4").Value = 0.2625365920782201
$ws.Range("M4").Value = 0.2659646019404605
$ws.Range("N4").Value = 1.767460534423432
$ws.Range("O4").Value = 3.367074478666083
# Row 5 (A5 index 3)
$ws.Range("B5").Value = 0.2451937588481883
$ws.Range("C5").Value = 0.03134856782351392
$ws.Range("D5").Value = 0.2527830071448562
$ws.Range("F5").Value = 1.467305972221922
$ws.Range("G5").Value = 0.002466074866024737
$ws.Range("I5").Value = 0.7060347426437446
$ws.Range("J5").Value = 0.2951664053675955
$ws.Range("K5").Value = 0.2541894231163155
$ws.Range("M5").Value = 0.2627886552157079
$ws.Range("N5").Value = 1.772057781195759
$ws.Range("O5").Value = 3.37215119906115
# Row 6 (A6 index 4)
$ws.Range("B6").Value = 0.2439277895241503
$ws.Range("C6").Value = 0.03114829255022755
$ws.Range("D6").Value = 0.2525875934061048
$ws.Range("F6").Value = 1.467507049676499
$ws.Range("G6").Value = 0.002466177122950834
$ws.Range("I6").Value = 0.7062922699632566
$ws.Range("J6").Value = 0.2951118818543321
$ws.Range("K6").Value = 0.2528038164563213
$ws.Range("M6").Value = 0.2622631265294615
$ws.Range("N6").Value = 1.772829254166869
$ws.Range("O6").Value = 3.37301854580096
# Row 7 (A7 index 5)
$ws.Range("B7").Value = 0.2527184712628525
$ws.Range("C7").Value = 0.03253820076760405
$ws.Range("D7").Value = 0.2539539573971865
$ws.Range("F7").Value = 1.466173712659696
$ws.Range("G7").Value = 0.00246547401243886
$ws.Range("I7").Value = 0.7045329852960052
$ws.Range("J7").Value = 0.2955035602929783
$ws.Range("K7").Value = 0.262423990437668
$ws.Range("M7").Value = 0.2659216471385051
$ws.Range("N7").Value = 1.767521991264334
$ws.Range("O7").Value = 3.367141312100358
# Row 8 (A8 index 6)
$ws.Range("B8").Value = 0.2915923995091703
$ws.Range("C8").Value = 0.03866572831489634
$ws.Range("D8").Value = 0.2602295293802683
$ws.Range("F8").Value = 1.461831936425469
$ws.Range("G8").Value = 0.002462535839075505
$ws.Range("I8").Value = 0.6974681397188505
$ws.Range("J8").Value = 0.2975589550049378
$ws.Range("K8").Value = 0.3049354176059182
$ws.Range("M8").Value = 0.2823311845116692
$ws.Range("N8").Value = 1.745279925052262
$ws.Range("O8").Value = 3.344745887940519
# Row 9 (A9 index 7)
$ws.Range("B9").Value = 0.3680927696585456
$ws.Range("C9").Value = 0.05065790139028081
$ws.Range("D9").Value = 0.2734133953940869
$ws.Range("F9").Value = 1.458849985112202
$ws.Range("G9").Value = 0.002457361654686311
$ws.Range("I9").Value = 0.6861225734236456
$ws.Range("J9").Value = 0.3027628655665566
$ws.Range("K9").Value = 0.3884864029271284
$ws.Range("M9").Value = 0.3154521397116383
$ws.Range("N9").Value = 1.7058937926708
$ws.Range("O9").Value = 3.313537691284637
# Row 10 (A10 index 8)
$ws.Range("B10").Value = 0.4244434241421686
$ws.Range("C10").Value = 0.05945247075663929
$ws.Range("D10").Value = 0.2836298196812095
$ws.Range("F10").Value = 1.460025502831655
$ws.Range("G10").Value = 0.002453914892542022
$ws.Range("I10").Value = 0.6793147768725376
$ws.Range("J10").Value = 0.307300709485574
$ws.Range("K10").Value = 0.4499671831358683
$ws.Range("M10").Value = 0.3403534367121281
$ws.Range("N10").Value = 1.679532289177632
$ws.Range("O10").Value = 3.298349262154233
# Row 11 (A11 index 9)
$ws.Range("B11").Value = 0.4501069686698997
$ws.Range("C11").Value = 0.06344958453964011
$ws.Range("D11").Value = 0.2883917471436206
$ws.Range("F11").Value = 1.461290059415603
$ws.Range("G11").Value = 0.002452423152650506
$ws.Range("I11").Value = 0.6765491855641628
$ws.Range("J11").Value = 0.3095201411469048
$ws.Range("K11").Value = 0.4779538049298537
$ws.Range("M11").Value = 0.3518035684714249
$ws.Range("N11").Value = 1.668098251177348
$ws.Range("O11").Value = 3.293118987511235
# Row 12 (A12 index 10)
$ws.Range("B12").Value = 0.4598288836229187
$ws.Range("C12").Value = 0.06496262978822642
$ws.Range("D12").Value = 0.2902113136351403
$ws.Range("F12").Value = 1.461873726047912
$ws.Range("G12").Value = 0.002451869172853973
$ws.Range("I12").Value = 0.6755495472749118
$ws.Range("J12").Value = 0.3103828651845504
$ws.Range("K12").Value = 0.4885538748607132
$ws.Range("M12").Value = 0.3561568653158886
$ws.Range("N12").Value = 1.663848688539018
$ws.Range("O12").Value = 3.291379706822141
# Row 13 (A13 index 11)
$ws.Range("B13").Value = 0.4577349397699493
$ws.Range("C13").Value = 0.06463679495355734
$ws.Range("D13").Value = 0.2898187132603454
$ws.Range("F13").Value = 1.461743363915588
$ws.Range("G13").Value = 0.002451987997890592
$ws.Range("I13").Value = 0.6757627190302955
$ws.Range("J13").Value = 0.310196072087038
$ws.Range("K13").Value = 0.4862708735067542
$ws.Range("M13").Value = 0.3552185354838855
$ws.Range("N13").Value = 1.664760337856888
$ws.Range("O13").Value = 3.291743561746671
# Row 14 (A14 index 12)
$ws.Range("B14").Value = 0.4509067256126684
$ws.Range("C14").Value = 0.0635740755528218
$ws.Range("D14").Value = 0.2885411174627421
$ws.Range("F14").Value = 1.461335978339747
$ws.Range("G14").Value = 0.00245237735816558
$ws.Range("I14").Value = 0.6764659901833632
$ws.Range("J14").Value = 0.3095906717873618
$ws.Range("K14").Value = 0.4788258392320017
$ws.Range("M14").Value = 0.3521613697638699
$ws.Range("N14").Value = 1.667747028414336
$ws.Range("O14").Value = 3.292971060359008
# Row 15 (A15 index 13)
$ws.Range("B15").Value = 0.4467247116965893
$ws.Range("C15").Value = 0.06292305256332043
$ws.Range("D15").Value = 0.2877606757707554
$ws.Range("F15").Value = 1.461100087680222
$ws.Range("G15").Value = 0.002452617271847562
$ws.Range("I15").Value = 0.676902966730335
$ws.Range("J15").Value = 0.3092227462074248
$ws.Range("K15").Value = 0.4742658037104945
$ws.Range("M15").Value = 0.3502910243486212
$ws.Range("N15").Value = 1.669586915636115
$ws.Range("O15").Value = 3.293754360633471
# Row 16 (A16 index 14)
$ws.Range("B16").Value = 0.4227667962003636
$ws.Range("C16").Value = 0.05919117278767771
$ws.Range("D16").Value = 0.2833209082581618
$ws.Range("F16").Value = 1.459957534812048
$ws.Range("G16").Value = 0.002454013910607171
$ws.Range("I16").Value = 0.6795021806384653
$ws.Range("J16").Value = 0.3071587832210554
$ws.Range("K16").Value = 0.4481385213165368
$ws.Range("M16").Value = 0.3396075878918055
$ws.Range("N16").Value = 1.680290760010564
$ws.Range("O16").Value = 3.298724847933272
# Row 17 (A17 index 15)
$ws.Range("B17").Value = 0.4080765291227237
$ws.Range("C17").Value = 0.05690082570265531
$ws.Range("D17").Value = 0.2806264719133083
$ws.Range("F17").Value = 1.459443438148057
$ws.Range("G17").Value = 0.002454890186986881
$ws.Range("I17").Value = 0.6811815627873976
$ws.Range("J17").Value = 0.3059323192567973
$ws.Range("K17").Value = 0.4321146824633217
$ws.Range("M17").Value = 0.3330848395902066
$ws.Range("N17").Value = 1.687000143340157
$ws.Range("O17").Value = 3.302204025283771
# Row 18 (A18 index 16)
$ws.Range("B18").Value = 0.3996298703961543
$ws.Range("C18").Value = 0.0555831452798401
$ws.Range("D18").Value = 0.2790874811662007
$ws.Range("F18").Value = 1.459216444626676
$ws.Range("G18").Value = 0.002455401374162481
$ws.Range("I18").Value = 0.6821786840349873
$ws.Range("J18").Value = 0.3052414946574658
$ws.Range("K18").Value = 0.4228999876044384
$ws.Range("M18").Value = 0.3293446620830025
$ws.Range("N18").Value = 1.690911717493687
$ws.Range("O18").Value = 3.30436320095788
# Row 19 (A19 index 17)
$ws.Range("B19").Value = 0.3967704730445405
$ws.Range("C19").Value = 0.05513694590935359
$ws.Range("D19").Value = 0.2785682602755628
$ws.Range("F19").Value = 1.45915139226129
$ws.Range("G19").Value = 0.002455575687170872
$ws.Range("I19").Value = 0.6825216485507397
$ws.Range("J19").Value = 0.3050101028399581
$ws.Range("K19").Value = 0.419780374918048
$ws.Range("M19").Value = 0.3280802900952722
$ws.Range("N19").Value = 1.692245127867074
$ws.Range("O19").Value = 3.305121409064157
# Row 20 (A20 index 18)
$ws.Range("B20").Value = 0.4096400477840518
$ws.Range("C20").Value = 0.05714467199703677
$ws.Range("D20").Value = 0.2809121847169536
$ws.Range("F20").Value = 1.459491055295288
$ws.Range("G20").Value = 0.002454796163262891
$ws.Range("I20").Value = 0.680999562232639
$ws.Range("J20").Value = 0.3060613671334522
$ws.Range("K20").Value = 0.4338202658059913
$ws.Range("M20").Value = 0.3337780050769794
$ws.Range("N20").Value = 1.686280482799535
$ws.Range("O20").Value = 3.301817304800977
# Row 21 (A21 index 19)
$ws.Range("B21").Value = 0.4529122417959854
$ws.Range("C21").Value = 0.06388623808652483
$ws.Range("D21").Value = 0.2889159361910458
$ws.Range("F21").Value = 1.461452794049919
$ws.Range("G21").Value = 0.002452262698046258
$ws.Range("I21").Value = 0.6762581297888559
$ws.Range("J21").Value = 0.3097678883246147
$ws.Range("K21").Value = 0.4810125724475824
$ws.Range("M21").Value = 0.3530588635108955
$ws.Range("N21").Value = 1.66686758654696
$ws.Range("O21").Value = 3.292603966192331
# Row 22 (A22 index 20)
$ws.Range("B22").Value = 0.4812143687522621
$ws.Range("C22").Value = 0.06828886174844229
$ws.Range("D22").Value = 0.2942419631940822
$ws.Range("F22").Value = 1.463345731428959
$ws.Range("G22").Value = 0.002450670497645357
$ws.Range("I22").Value = 0.6734369552631172
$ws.Range("J22").Value = 0.3123201238779103
$ws.Range("K22").Value = 0.511867713824671
$ws.Range("M22").Value = 0.3657612223804918
$ws.Range("N22").Value = 1.654647964581859
$ws.Range("O22").Value = 3.287989000419088
# Row 23 (A23 index 21)
$ws.Range("B23").Value = 0.4661072230150864
$ws.Range("C23").Value = 0.06593942777712414
$ws.Range("D23").Value = 0.2913907013195995
$ws.Range("F23").Value = 1.462279590370841
$ws.Range("G23").Value = 0.002451514484732279
$ws.Range("I23").Value = 0.6749172699737187
$ws.Range("J23").Value = 0.3109460821205943
$ws.Range("K23").Value = 0.4953988047887208
$ws.Range("M23").Value = 0.3589725476414714
$ws.Range("N23").Value = 1.661126990617786
$ws.Range("O23").Value = 3.290323444146139
# Row 24 (A24 index 22)
$ws.Range("B24").Value = 0.4089331841610431
$ws.Range("C24").Value = 0.05703443197752733
$ws.Range("D24").Value = 0.2807829825615187
$ws.Range("F24").Value = 1.459469313989544
$ws.Range("G24").Value = 0.002454838648437199
$ws.Range("I24").Value = 0.6810817461815866
$ws.Range("J24").Value = 0.3060029800873849
$ws.Range("K24").Value = 0.4330491789378357
$ws.Range("M24").Value = 0.3334645943019865
$ws.Range("N24").Value = 1.6866056723127
$ws.Range("O24").Value = 3.301991646000005
# Row 25 (A25 index 23)
$ws.Range("B25").Value = 0.3473704828043083
$ws.Range("C25").Value = 0.04741641493393445
$ws.Range("D25").Value = 0.2697533542298487
$ws.Range("F25").Value = 1.45906506174709
$ws.Range("G25").Value = 0.002458698867501426
$ws.Range("I25").Value = 0.6889234457483298
$ws.Range("J25").Value = 0.3012295366269839
$ws.Range("K25").Value = 0.3658655355862663
$ws.Range("M25").Value = 0.306391915543081
$ws.Range("N25").Value = 1.716096257110493
$ws.Range("O25").Value = 3.320620520761622
